$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "...Mock data provide" paragraph: append "d of courses" + "." as two new
#    runs (matching Word's behaviour of creating fresh, rPr-less runs when
#    text is typed at the end of a paragraph), then re-create the hidden
#    "_GoBack" bookmark at the very end of the paragraph (Word moves this
#    bookmark to the most recent edit location; re-adding it under the same
#    name automatically removes it from its previous location later in the
#    document).
# ---------------------------------------------------------------------------

$para = $d.Paragraphs(5)
$rng = $para.Range
$rng.End = $rng.End - 1
$rng.InsertAfter("d of courses")
$rng.InsertAfter(".")

# Placing a bookmark collapsed exactly on the character right before a
# paragraph mark is unreliable, so append a throw-away sentinel character,
# anchor the bookmark just before it (a safe, non-boundary position), and
# then delete the sentinel again.
$rng.InsertAfter("X")

$para3 = $d.Paragraphs(5)
$sentinelPos = $para3.Range.End - 2
$bmRng = $d.Range($sentinelPos, $sentinelPos)
$d.Bookmarks.Add("_GoBack", $bmRng)

$para4 = $d.Paragraphs(5)
$sentinelPos2 = $para4.Range.End - 2
$sentinelRng = $d.Range($sentinelPos2, $sentinelPos2 + 1)
$sentinelRng.Text = ""

# ---------------------------------------------------------------------------
# 2) Table cells: merge the split "[" / "Yingluck " and "[" / "Mark " runs
#    (and the ", " / "Mark " pair) back into single runs. The visible text is
#    unchanged - only the run boundary goes away - so search/replace with the
#    exact same text forces Word to normalise the (identically-formatted)
#    adjacent runs into one. The search/replace strings intentionally avoid
#    leading/trailing whitespace so the writer does not tag the merged run
#    with xml:space="preserve", and they avoid the non-breaking space so that
#    character is left completely untouched.
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("[Yingluck", $true, $false, $false, $false, $false, $true, 1, $false, "[Yingluck", 2) | Out-Null
$d.Content.Find.Execute(", Mark", $true, $false, $false, $false, $false, $true, 1, $false, ", Mark", 2) | Out-Null
$d.Content.Find.Execute("[Mark", $true, $false, $false, $false, $false, $true, 1, $false, "[Mark", 2) | Out-Null
